# Apply updated "dSF" (column F) values on Sheet1, per repulled data.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Map of row number -> new value for column F (dSF)
$updates = @{
    7  = 0
    9  = 0
    10 = 0
    14 = -5
    16 = 1
    25 = -1
    28 = 4
    40 = 0
    44 = -1
    45 = -2
    50 = 3
    54 = -1
    58 = -3
    60 = -4
    61 = -5
    70 = -2
    77 = -1
    79 = -1
    81 = 0
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
